# "Update country data files" - refine precision of several MSME indicator
# values for Bosnia and Herzegovina (Summary sheet).
#
# These cells hold numeric-looking values that are stored as TEXT (shared
# strings) in the workbook, e.g. "39.1" rather than the number 39.1. Set the
# cell's NumberFormat to Text ("@") before writing so Excel keeps the new
# value as text too, instead of silently converting it to a numeric cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "39.13"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "2.64"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.77"

# Employment (% of total): MSMEs
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.83"

# Enterprises (% of total): Micro / SMEs / MSMEs
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "93.27"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "6.29"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "99.56"
